# daily auto push: 2026-01-31 09:39 UTC
# Insert a new data row at row 756 (pushing existing rows 756:797 down to
# 757:798) and populate it with the new day's first observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 756 and below down by one row, then fill in the new row's data.
$ws.Rows.Item(756).Insert()

# Force the date-like text into column A as a literal string (matching the
# other "YYYY/MM/DD" text cells in the sheet) instead of letting Excel
# auto-convert it to a date serial number, then restore the default
# (unstyled) cell format used by every other data row.
$ws.Cells.Item(756, 1).NumberFormat = "@"
$ws.Cells.Item(756, 1).Value = "2026/01/31"
$ws.Cells.Item(756, 1).Style = "Normal"

$ws.Cells.Item(756, 2).Value = "土"
$ws.Cells.Item(756, 3).Value = 16
$ws.Cells.Item(756, 4).Value = 188
